$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the MobileNumberField_Xpath locator value (row 13, col B) ---
$ws.Cells.Item(13, 2).Value = "//input[@name='socialogin_email']"

# --- Append new rows 16-24: "Review and Ratings" test case locators ---
$data = @(
    @("Login_For_Review_And_Rating_LinkText", "Log in"),
    @("Rating_Text_Xpath", "(//form[@id='review-form']//div[2])[1]/h5"),
    @("Review_Text_Xpath", "(//form[@id='review-form']//div[2])[1]/ul/h5"),
    @("Rating_Stars_Xpath", "(//form[@id='review-form']//div[2])[1]/div[1]/ul/li"),
    @("Empty_Rating_Stars_Xpath", "(//form[@id='review-form']//div[2])[1]/div[1]/div[2]/input"),
    @("Review_Title_Label_Xpath", "//*[@id='review-form']/fieldset/div[2]/ul/li[1]/label"),
    @("Review_Title_TextBox_Xpath", "//*[@id='summary_field']"),
    @("Review_Xpath", "//*[@id='review-form']/fieldset/div[2]/ul/li[2]/label"),
    @("Review_TextBox_Xpath", "//*[@id='review_field']")
)

$row = 16
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}

# Row 16's key cell (Login_For_Review_And_Rating_LinkText) holds a link text
# rather than an xpath, so - mirroring the rest of the sheet - it picks up
# the blue "Courier New" styling normally reserved for the xpath/value
# column. Copy that look from an existing xpath-value cell.
$ws.Range("B2").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B25").Select()
